$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating-point precision on the existing A2 timestamp
$ws.Range("A2").Value = 45863.37517840278

# Append a new data row (row 3) written by the scheduled task run
$ws.Range("A3").Value = 45863.41686092952
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("B3").Value = 2025
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 16.5
$ws.Range("E3").Value = 82.47
$ws.Range("F3").Value = 462.3
$ws.Range("G3").Value = 8.539999999999999
$ws.Range("H3").Value = "ESE"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "10:00:16"
